$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings stay stored as text (matches source data's inline-string cells)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "30.629.82"
$ws.Range("E2").Value = "  +0.83%  "
$ws.Range("D3").Value = "2.114.84"
$ws.Range("E3").Value = "  +0.52%  "
$ws.Range("E4").Value = "  +0.63%  "
$ws.Range("D5").Value = "349.85"
$ws.Range("E5").Value = "  +4.58%  "
$ws.Range("E6").Value = "  +0.80%  "
$ws.Range("D7").Value = "0.5251"
$ws.Range("E7").Value = "  +0.43%  "
$ws.Range("D8").Value = "0.4511"
$ws.Range("E8").Value = "  -1.25%  "
$ws.Range("D9").Value = "53.81"
$ws.Range("E9").Value = "  +1.45%  "
$ws.Range("D10").Value = "0.09024"
$ws.Range("E10").Value = "  +1.28%  "
$ws.Range("D11").Value = "1.173"
$ws.Range("E11").Value = "  -0.24%  "
$ws.Range("D12").Value = "24.50"
$ws.Range("E12").Value = "  +0.40%  "
$ws.Range("D13").Value = "2.125.81"
$ws.Range("E13").Value = "  +1.51%  "
$ws.Range("D14").Value = "6.830"
$ws.Range("E14").Value = "  +0.48%  "
$ws.Range("D15").Value = "8.031"
$ws.Range("E15").Value = "  +0.67%  "
$ws.Range("E16").Value = "  +5.81%  "
$ws.Range("D17").Value = "0.00001169"
$ws.Range("E17").Value = "  +3.26%  "
$ws.Range("E18").Value = "  +0.79%  "
$ws.Range("D19").Value = "0.06710"
$ws.Range("E19").Value = "  +1.19%  "
$ws.Range("D20").Value = "19.41"
$ws.Range("E20").Value = "  +0.59%  "
$ws.Range("D21").Value = "1.009"
$ws.Range("E21").Value = "  +0.75%  "
$ws.Range("D22").Value = "6.305"
$ws.Range("E22").Value = "  -0.56%  "
$ws.Range("D23").Value = "30.705.66"
$ws.Range("E23").Value = "  +0.83%  "
$ws.Range("D24").Value = "12.83"
$ws.Range("E24").Value = "  +3.46%  "
$ws.Range("D25").Value = "2.389"
$ws.Range("E25").Value = "  +1.18%  "
$ws.Range("D26").Value = "2.373.79"
$ws.Range("E26").Value = "  +1.31%  "
$ws.Range("D27").Value = "22.37"
$ws.Range("E27").Value = "  +0.28%  "
$ws.Range("D28").Value = "165.28"
$ws.Range("E28").Value = "  +0.95%  "
$ws.Range("D29").Value = "2.539"
$ws.Range("E29").Value = "  -1.10%  "
$ws.Range("D30").Value = "136.26"
$ws.Range("E30").Value = "  +2.77%  "
$ws.Range("D31").Value = "1.192"
$ws.Range("E31").Value = "  -2.60%  "
$ws.Range("D32").Value = "0.1077"
$ws.Range("D33").Value = "1.644"
$ws.Range("E33").Value = "  -4.95%  "
$ws.Range("D34").Value = "6.384"
$ws.Range("E34").Value = "  +2.98%  "
$ws.Range("E35").Value = "  +2.36%  "
$ws.Range("D36").Value = "10.33"
$ws.Range("E36").Value = "  -1.67%  "
$ws.Range("D37").Value = "5.925"
$ws.Range("E37").Value = "  +6.92%  "
$ws.Range("D38").Value = "0.02657"
$ws.Range("D39").Value = "0.06843"
$ws.Range("E39").Value = "  +0.22%  "
$ws.Range("D40").Value = "0.2318"
$ws.Range("E40").Value = "  +0.82%  "
$ws.Range("D41").Value = "12.54"
$ws.Range("E41").Value = "  -1.96%  "
$ws.Range("D42").Value = "0.6890"
$ws.Range("E42").Value = "  +0.01%  "
$ws.Range("D43").Value = "1.272"
$ws.Range("E43").Value = "  +2.17%  "
$ws.Range("D44").Value = "14.71"
$ws.Range("E44").Value = "  +4.60%  "
$ws.Range("D45").Value = "2.329"
$ws.Range("E45").Value = "  -0.94%  "
$ws.Range("D46").Value = "0.6447"
$ws.Range("E46").Value = "  +1.24%  "
$ws.Range("D48").Value = "0.00000000356"
$ws.Range("E48").Value = "  +1.01%  "
$ws.Range("D49").Value = "1.251"
$ws.Range("E49").Value = "  +0.44%  "
$ws.Range("D50").Value = "0.07295"
$ws.Range("E50").Value = "  +2.42%  "
$ws.Range("D51").Value = "82.47"
$ws.Range("E51").Value = "  -1.21%  "
